$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H11").Value = 1345.25
$ws.Range("I11").Value = 1345.25
$ws.Range("K11").Value = 1345.25
$ws.Range("M11").Value = -1205.25

$ws.Range("H32").Value = 3945.9023
$ws.Range("I32").Value = 3684.4443
$ws.Range("J32").Value = 4019.4375
$ws.Range("K32").Value = 3684.4443
$ws.Range("L32").Value = 4019.4375
$ws.Range("M32").Value = -3358.4443
$ws.Range("N32").Value = -4671.4375

$ws.Range("H40").Value = 6151.553
$ws.Range("I40").Value = 11288.272
$ws.Range("J40").Value = 4582
$ws.Range("K40").Value = 11288.272
$ws.Range("L40").Value = 4582
$ws.Range("M40").Value = -11113.272
$ws.Range("N40").Value = -4932

$ws.Range("H86").Value = 5541.125
$ws.Range("I86").Value = 5146.875
$ws.Range("J86").Value = 5738.25
$ws.Range("K86").Value = 5146.875
$ws.Range("L86").Value = 5738.25
$ws.Range("M86").Value = -4023.875
$ws.Range("N86").Value = -7984.25

$ws.Range("H89").Value = 5541.125
$ws.Range("I89").Value = 5146.875
$ws.Range("J89").Value = 5738.25
$ws.Range("K89").Value = 25734.375
$ws.Range("L89").Value = 28691.25
$ws.Range("M89").Value = -20118.375
$ws.Range("N89").Value = -39923.25

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H118").Value = 332.64285
$ws.Range("I118").Value = 332.64285
$ws.Range("K118").Value = 997.9285500000001
$ws.Range("M118").Value = 659.0714499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3581.3076
$ws.Range("I32").Value = 2252.6182
$ws.Range("K32").Value = 2252.6182
$ws.Range("M32").Value = -1965.6182

$ws.Range("H34").Value = 20000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20542
$ws.Range("M34").ClearContents()

$ws.Range("H45").Value = 7995989.5
$ws.Range("I45").Value = 15985291
$ws.Range("J45").Value = 6688.1113
$ws.Range("K45").Value = 15985291
$ws.Range("L45").Value = 6688.1113
$ws.Range("M45").Value = -15984914
$ws.Range("N45").Value = -7442.1113

$ws.Range("H61").Value = 15998.75
$ws.Range("I61").Value = 22599
$ws.Range("J61").Value = 4998.3335
$ws.Range("K61").Value = 22599
$ws.Range("L61").Value = 4998.3335
$ws.Range("M61").Value = -22387
$ws.Range("N61").Value = -5422.3335

$ws.Range("H74").Value = 57008.55
$ws.Range("I74").Value = 11888.852
$ws.Range("K74").Value = 11888.852
$ws.Range("M74").Value = -11014.852

$ws.Range("H77").Value = 57008.55
$ws.Range("I77").Value = 11888.852
$ws.Range("K77").Value = 59444.26
$ws.Range("M77").Value = -55076.26

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H102").Value = 5213143.5
$ws.Range("I102").Value = 5957307.5
$ws.Range("K102").Value = 5957307.5
$ws.Range("M102").Value = -5955685.5

$ws.Range("H132").Value = 18335
$ws.Range("J132").Value = 13497.5
$ws.Range("L132").Value = 40492.5
$ws.Range("N132").Value = -45552.5

$ws.Range("H136").Value = 15998.75
$ws.Range("I136").Value = 22599
$ws.Range("J136").Value = 4998.3335
$ws.Range("K136").Value = 67797
$ws.Range("L136").Value = 14995.0005
$ws.Range("M136").Value = -65247
$ws.Range("N136").Value = -20095.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13891306
$ws.Range("I20").Value = 17545644
$ws.Range("K20").Value = 17545644
$ws.Range("M20").Value = -17545397

$ws.Range("H22").Value = 819.86365
$ws.Range("J22").Value = 281
$ws.Range("L22").Value = 281
$ws.Range("N22").Value = -627

$ws.Range("H86").Value = 6264937.5
$ws.Range("I86").Value = 6682500
$ws.Range("J86").Value = 1497
$ws.Range("K86").Value = 6682500
$ws.Range("L86").Value = 1497
$ws.Range("M86").Value = -6681377
$ws.Range("N86").Value = -3743

$ws.Range("H89").Value = 6264937.5
$ws.Range("I89").Value = 6682500
$ws.Range("J89").Value = 1497
$ws.Range("K89").Value = 33412500
$ws.Range("L89").Value = 7485
$ws.Range("M89").Value = -33406884
$ws.Range("N89").Value = -18717

$ws.Range("H105").Value = 3175880.2
$ws.Range("I105").Value = 3342979.5
$ws.Range("K105").Value = 3342979.5
$ws.Range("M105").Value = -3341232.5

$ws.Range("H134").Value = 29998.666
$ws.Range("I134").Value = 30713.285
$ws.Range("K134").Value = 92139.855
$ws.Range("M134").Value = -89604.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17243.363
$ws.Range("I31").Value = 8878.691999999999
$ws.Range("J31").Value = 18866.357
$ws.Range("K31").Value = 8878.691999999999
$ws.Range("L31").Value = 18866.357
$ws.Range("M31").Value = -8583.691999999999
$ws.Range("N31").Value = -19456.357

$ws.Range("H32").Value = 1282
$ws.Range("I32").Value = 1282
$ws.Range("K32").Value = 1282
$ws.Range("M32").Value = -966

$ws.Range("H34").Value = 17243.363
$ws.Range("I34").Value = 8878.691999999999
$ws.Range("J34").Value = 18866.357
$ws.Range("K34").Value = 8878.691999999999
$ws.Range("L34").Value = 18866.357
$ws.Range("M34").Value = -8676.691999999999
$ws.Range("N34").Value = -19270.357

$ws.Range("H105").Value = 1252.0834
$ws.Range("I105").Value = 1126.6666
$ws.Range("J105").Value = 2130
$ws.Range("K105").Value = 1126.6666
$ws.Range("L105").Value = 2130
$ws.Range("M105").Value = 620.3334
$ws.Range("N105").Value = -5624

$ws.Range("H132").Value = 46282.043
$ws.Range("I132").Value = 47340.363
$ws.Range("K132").Value = 142021.089
$ws.Range("M132").Value = -139491.089

$ws.Range("H134").Value = 9391.585999999999
$ws.Range("I134").Value = 6961.1055
$ws.Range("J134").Value = 14009.5
$ws.Range("K134").Value = 20883.3165
$ws.Range("L134").Value = 42028.5
$ws.Range("M134").Value = -18348.3165
$ws.Range("N134").Value = -47098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5761.15
$ws.Range("I5").Value = 830.7
$ws.Range("K5").Value = 2492.1
$ws.Range("M5").Value = -2380.1

$ws.Range("H11").Value = 83992.16
$ws.Range("I11").Value = 4862.25
$ws.Range("K11").Value = 14586.75
$ws.Range("M11").Value = -14446.75

$ws.Range("H135").Value = 5761.15
$ws.Range("I135").Value = 830.7
$ws.Range("K135").Value = 7476.3
$ws.Range("M135").Value = -4941.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10004585
$ws.Range("I70").Value = 16671142
$ws.Range("K70").Value = 16671142
$ws.Range("M70").Value = -16670872

$ws.Range("H73").Value = 10004585
$ws.Range("I73").Value = 16671142
$ws.Range("K73").Value = 16671142
$ws.Range("M73").Value = -16670206

$ws.Range("H80").Value = 43569292
$ws.Range("I80").Value = 59220990
$ws.Range("K80").Value = 59220990
$ws.Range("M80").Value = -59219992

$ws.Range("H83").Value = 43569292
$ws.Range("I83").Value = 59220990
$ws.Range("K83").Value = 296104950
$ws.Range("M83").Value = -296099958

$ws.Range("H102").Value = 5230430.5
$ws.Range("I102").Value = 8548925
$ws.Range("J102").Value = 1911936.5
$ws.Range("K102").Value = 8548925
$ws.Range("L102").Value = 1911936.5
$ws.Range("M102").Value = -8547303
$ws.Range("N102").Value = -1915180.5

$ws.Range("H122").Value = 282093.62
$ws.Range("I122").Value = 359862.62
$ws.Range("K122").Value = 1079587.86
$ws.Range("M122").Value = -1077137.86

$ws.Range("H132").Value = 12116.479
$ws.Range("J132").Value = 35665
$ws.Range("L132").Value = 106995
$ws.Range("N132").Value = -112055

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4489.952
$ws.Range("I100").Value = 3987.647
$ws.Range("J100").Value = 6624.75
$ws.Range("K100").Value = 3987.647
$ws.Range("L100").Value = 6624.75
$ws.Range("M100").Value = -3446.647
$ws.Range("N100").Value = -7706.75

$ws.Range("H136").Value = 28097.375
$ws.Range("I136").Value = 39496.965
$ws.Range("J136").Value = 4421.3076
$ws.Range("K136").Value = 118490.895
$ws.Range("L136").Value = 13263.9228
$ws.Range("M136").Value = -115940.895
$ws.Range("N136").Value = -18363.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 27782022
$ws.Range("I81").Value = 55556990
$ws.Range("J81").Value = 7053
$ws.Range("K81").Value = 111113980
$ws.Range("L81").Value = 14106
$ws.Range("M81").Value = -111112919
$ws.Range("N81").Value = -16228

$ws.Range("H84").Value = 27782022
$ws.Range("I84").Value = 55556990
$ws.Range("J84").Value = 7053
$ws.Range("K84").Value = 555569900
$ws.Range("L84").Value = 70530
$ws.Range("M84").Value = -555564596
$ws.Range("N84").Value = -81138

$ws.Range("H122").Value = 2826.6
$ws.Range("I122").Value = 2768.0435
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8304.130500000001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -5854.130500000001
$ws.Range("N122").Value = -15400

$ws.Range("H126").Value = 2368.926
$ws.Range("J126").Value = 3791.8
$ws.Range("L126").Value = 11375.4
$ws.Range("N126").Value = -16315.4

$ws.Range("H132").Value = 18386758
$ws.Range("I132").Value = 26321686
$ws.Range("J132").Value = 649863.0600000001
$ws.Range("K132").Value = 78965058
$ws.Range("L132").Value = 1949589.18
$ws.Range("M132").Value = -78962528
$ws.Range("N132").Value = -1954649.18

$ws.Range("H136").Value = 4990.849
$ws.Range("I136").Value = 6734.9653
$ws.Range("J136").Value = 2883.375
$ws.Range("K136").Value = 20204.8959
$ws.Range("L136").Value = 8650.125
$ws.Range("M136").Value = -17654.8959
$ws.Range("N136").Value = -13750.125
